$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.069.93"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.301.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.520"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.50%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.657.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.299.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.789"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.941.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "238.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.04%  "
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -13.16%  "
$ws.Range("E30").Value = "  +0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "163.47"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0698"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.77"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.011.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.529.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("E51").Value = "  -0.21%  "
